$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.880.63'
$ws.Range('E2').Value = '  +0.37%  '
$ws.Range('D3').Value = '2.221.95'
$ws.Range('E3').Value = '  +0.46%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'291.93"
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').Value = "'87.05"
$ws.Range('E6').Value = '  +0.49%  '
$ws.Range('E7').Value = '  -0.43%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = "'0.467"
$ws.Range('E9').Value = '  -0.72%  '
$ws.Range('D10').Value = "'30.47"
$ws.Range('E10').Value = '  +0.78%  '
$ws.Range('D11').Value = "'0.0780"
$ws.Range('E11').Value = '  -0.60%  '
$ws.Range('E12').Value = '  +5.44%  '
$ws.Range('E13').Value = '  +2.70%  '
$ws.Range('D14').Value = "'6.44"
$ws.Range('E14').Value = '  +1.65%  '
$ws.Range('D15').Value = '2.569.17'
$ws.Range('E15').Value = '  +0.65%  '
$ws.Range('D16').Value = "'13.81"
$ws.Range('E16').Value = '  -1.53%  '
$ws.Range('D17').Value = '2.221.85'
$ws.Range('E17').Value = '  +1.18%  '
$ws.Range('D18').Value = "'0.731"
$ws.Range('E18').Value = '  +0.63%  '
$ws.Range('D19').Value = '39.835.03'
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('D20').Value = '0.0₃0885'
$ws.Range('E20').Value = '  +0.62%  '
$ws.Range('D21').Value = "'11.09"
$ws.Range('E21').Value = '  -3.44%  '
$ws.Range('E22').Value = '  -0.76%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').Value = "'237.03"
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('E26').Value = '  -0.31%  '
$ws.Range('E27').Value = '  -0.37%  '
$ws.Range('D28').Value = "'23.05"
$ws.Range('E28').Value = '  +1.43%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = "'2.14"
$ws.Range('E29').Value = '  -2.96%  '
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').Value = "'9.25"
$ws.Range('E30').Value = '  -0.21%  '
$ws.Range('D31').Value = "'156.24"
$ws.Range('E31').Value = '  +2.59%  '
$ws.Range('D32').Value = "'31.93"
$ws.Range('E32').Value = '  -2.40%  '
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('E34').Value = '  +0.93%  '
$ws.Range('E35').Value = '  +7.23%  '
$ws.Range('D36').Value = "'0.0714"
$ws.Range('E36').Value = '  -0.28%  '
$ws.Range('E37').Value = '  -1.70%  '
$ws.Range('E38').Value = '  -0.27%  '
$ws.Range('E39').Value = '  +3.30%  '
$ws.Range('D40').Value = "'0.0994"
$ws.Range('E40').Value = '  +0.84%  '
$ws.Range('E41').Value = '  -4.32%  '
$ws.Range('D42').Value = '2.109.66'
$ws.Range('E42').Value = '  +1.41%  '
$ws.Range('D43').Value = "'3.71"
$ws.Range('E43').Value = '  -1.76%  '
$ws.Range('D44').Value = "'18.16"
$ws.Range('E44').Value = '  +2.65%  '
$ws.Range('E45').Value = '  +1.24%  '
$ws.Range('D46').Value = "'9.89"
$ws.Range('E46').Value = '  -1.20%  '
$ws.Range('D47').Value = "'1.99"
$ws.Range('E47').Value = '  -7.98%  '
$ws.Range('D48').Value = "'2.72"
$ws.Range('E48').Value = '  +4.59%  '
$ws.Range('D49').Value = '2.439.72'
$ws.Range('E49').Value = '  +0.63%  '
$ws.Range('E51').Value = '  +2.37%  '
